$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ColumnBlock($col, $startRow, [object[]]$values) {
    $n = $values.Count
    $arr = New-Object "object[,]" $n,1
    for ($i = 0; $i -lt $n; $i++) { $arr[$i,0] = $values[$i] }
    $endRow = $startRow + $n - 1
    $rng = $col + $startRow + ":" + $col + $endRow
    $ws.Range($rng).Value = $arr
}

# LSNS-R block (incl. complete row): columns A-D as bulk column writes
Set-ColumnBlock "A" 399 @('lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr', 'lubben_social_network_scale_revised_lsnsr')
Set-ColumnBlock "B" 399 @('lsns_r_q1', 'lsns_r_q2', 'lsns_r_q3', 'lsns_r_q4', 'lsns_r_q5', 'lsns_r_q6', 'lsns_r_q7', 'lsns_r_q8', 'lsns_r_q9', 'lsns_r_q10', 'lsns_r_q11', 'lsns_r_q12', 'lubben_social_network_scale_revised_lsnsr_complete')
Set-ColumnBlock "C" 399 @('LUBBEN SOCIAL NETWORK SCALE - REVISED (LSNS-R) FAMILY: Considering the people to whom you are related by birth, marriage, adoption, etc…', 'LUBBEN SOCIAL NETWORK SCALE - REVISED (LSNS-R) FAMILY: Considering the people to whom you are related by birth, marriage, adoption, etc…', 'LUBBEN SOCIAL NETWORK SCALE - REVISED (LSNS-R) FAMILY: Considering the people to whom you are related by birth, marriage, adoption, etc…', 'LUBBEN SOCIAL NETWORK SCALE - REVISED (LSNS-R) FAMILY: Considering the people to whom you are related by birth, marriage, adoption, etc…', 'LUBBEN SOCIAL NETWORK SCALE - REVISED (LSNS-R) FAMILY: Considering the people to whom you are related by birth, marriage, adoption, etc…', 'LUBBEN SOCIAL NETWORK SCALE - REVISED (LSNS-R) FAMILY: Considering the people to whom you are related by birth, marriage, adoption, etc…', 'FRIENDSHIPS: Considering all of your friends including those who live in your neighbourhood...', 'FRIENDSHIPS: Considering all of your friends including those who live in your neighbourhood...', 'FRIENDSHIPS: Considering all of your friends including those who live in your neighbourhood...', 'FRIENDSHIPS: Considering all of your friends including those who live in your neighbourhood...', 'FRIENDSHIPS: Considering all of your friends including those who live in your neighbourhood...', 'FRIENDSHIPS: Considering all of your friends including those who live in your neighbourhood...', 'Form Status')
Set-ColumnBlock "D" 399 @('1. How many relatives do you see or hear from at least once a month?', '2. How often do you see or hear from the relative with whom you have the most contact?', '3. How many relatives do you feel at ease with that you can talk about private matters?', '4. How many relatives do you feel close to such that you could call on them for help?', '5. When one of your relatives has an important decision to make, how often do they talk to you about it?', '6. How often is one of your relatives available for you to talk to when you have an important decision to make?', '7. How many of your friends do you see or hear from at least once a month?', '8. How often do you see or hear from the friend with whom you have the most contact?', '9. How many friends do you feel at ease with that you can talk about private matters?', '10. How many friends do you feel close to such that you could call on them for help?', '11. When one of your friends has an important decision to make, how often do they talk to you about it?', '12. How often is one of your friends available for you to talk to when you have an important decision to make?', 'Complete?')

# LSNS-R block (incl. complete row): columns E,F row-major (per-row)
$ws.Range("E399").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F399").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E400").Value = '0 = less than monthly; 1 = monthly; 2 = few times a month; 3 = weekly; 4 = few times a week; 5 = daily'
$ws.Range("F400").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E401").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F401").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E402").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F402").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E403").Value = '0 = never; 1 = seldom; 2 = sometimes; 3 = often; 4 = very often; 5 = always'
$ws.Range("F403").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E404").Value = '0 = never; 1 = seldom; 2 = sometimes; 3 = often; 4 = very often; 5 = always'
$ws.Range("F404").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E405").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F405").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E406").Value = '0 = less than monthly; 1 = monthly; 2 = few times a month; 3 = weekly; 4 = few times a week; 5 = daily'
$ws.Range("F406").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E407").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F407").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E408").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F408").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E409").Value = '0 = never; 1 = seldom; 2 = sometimes; 3 = often; 4 = very often; 5 = always'
$ws.Range("F409").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E410").Value = '0 = never; 1 = seldom; 2 = sometimes; 3 = often; 4 = very often; 5 = always'
$ws.Range("F410").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E411").Value = 'Incomplete; Unverified; Complete'
$ws.Range("F411").Value = '0; 1; 2'

# LSNS-6 block (questions only): columns A-D as bulk column writes
Set-ColumnBlock "A" 412 @('lubben_social_network_scale_6_lsns_6', 'lubben_social_network_scale_6_lsns_6', 'lubben_social_network_scale_6_lsns_6', 'lubben_social_network_scale_6_lsns_6', 'lubben_social_network_scale_6_lsns_6', 'lubben_social_network_scale_6_lsns_6')
Set-ColumnBlock "B" 412 @('lsns_6_q1', 'lsns_6_q2', 'lsns_6_q3', 'lsns_6_q4', 'lsns_6_q5', 'lsns_6_q6')
Set-ColumnBlock "C" 412 @('LUBBEN SOCIAL NETWORK SCALE - REVISED (LSNS-R) FAMILY: Considering the people to whom you are related by birth, marriage, adoption, etc...', 'LUBBEN SOCIAL NETWORK SCALE - REVISED (LSNS-R) FAMILY: Considering the people to whom you are related by birth, marriage, adoption, etc...', 'LUBBEN SOCIAL NETWORK SCALE - REVISED (LSNS-R) FAMILY: Considering the people to whom you are related by birth, marriage, adoption, etc...', 'FRIENDSHIPS: Considering all of your friends including those who live in your neighbourhood...', 'FRIENDSHIPS: Considering all of your friends including those who live in your neighbourhood...', 'FRIENDSHIPS: Considering all of your friends including those who live in your neighbourhood...')
Set-ColumnBlock "D" 412 @('1. How many relatives do you see or hear from at least once a month?', '2. How many relatives do you feel at ease with that you can talk about private matters?', '3. How many relatives do you feel close to such that you could call on them for help?', '4. How many of your friends do you see or hear from at least once a month?', '5. How many friends do you feel at ease with that you can talk about private matters?', '6. How many friends do you feel close to such that you could call on them for help?')

# LSNS-6 block (questions only): columns E,F row-major (per-row)
$ws.Range("E412").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F412").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E413").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F413").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E414").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F414").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E415").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F415").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E416").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F416").Value = '0; 1; 2; 3; 4; 5'
$ws.Range("E417").Value = '0 = none; 1 = one; 2 = two; 3 = three or four; 4 = five thru eight; 5 = nine or more'
$ws.Range("F417").Value = '0; 1; 2; 3; 4; 5'

# LSNS-6 complete row (added separately): single row, all columns
$ws.Range("A418").Value = 'lubben_social_network_scale_6_lsns_6'
$ws.Range("B418").Value = 'lubben_social_network_scale_6_lsns_6_complete'
$ws.Range("C418").Value = 'Form Status'
$ws.Range("D418").Value = 'Complete?'
$ws.Range("E418").Value = 'Incomplete; Unverified; Complete'
$ws.Range("F418").Value = '0; 1; 2'

# Row heights to match wrapped-text auto-fit appearance
$ws.Rows("399").RowHeight = 43.2
$ws.Rows("400").RowHeight = 57.6
$ws.Rows("401").RowHeight = 43.2
$ws.Rows("402").RowHeight = 43.2
$ws.Rows("403").RowHeight = 43.2
$ws.Rows("404").RowHeight = 43.2
$ws.Rows("405").RowHeight = 43.2
$ws.Rows("406").RowHeight = 57.6
$ws.Rows("407").RowHeight = 43.2
$ws.Rows("408").RowHeight = 43.2
$ws.Rows("409").RowHeight = 43.2
$ws.Rows("410").RowHeight = 43.2
$ws.Rows("412").RowHeight = 43.2
$ws.Rows("413").RowHeight = 43.2
$ws.Rows("414").RowHeight = 43.2
$ws.Rows("415").RowHeight = 43.2
$ws.Rows("416").RowHeight = 43.2
$ws.Rows("417").RowHeight = 43.2
